$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "58.858.75"
$c.Style = "Normal"

$ws.Range("E2").Value = "  -0.37%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.498.31"
$c.Style = "Normal"

$ws.Range("E3").Value = "  -0.08%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"

$ws.Range("E4").Value = "  +0.20%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "537.03"
$c.Style = "Normal"

$ws.Range("E5").Value = "  -0.36%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "136.30"
$c.Style = "Normal"

$ws.Range("E6").Value = "  -1.45%  "

$ws.Range("E7").Value = "  -0.16%  "

$ws.Range("E8").Value = "  +0.99%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.522.17"
$c.Style = "Normal"

$ws.Range("E9").Value = "  +0.91%  "

$ws.Range("E10").Value = "  +0.90%  "

$ws.Range("E11").Value = "  -2.47%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "5.30"
$c.Style = "Normal"

$ws.Range("E12").Value = "  -2.02%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.347"
$c.Style = "Normal"

$ws.Range("E13").Value = "  -0.43%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "2.967.53"
$c.Style = "Normal"

$ws.Range("E14").Value = "  +0.65%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "22.99"
$c.Style = "Normal"

$ws.Range("E15").Value = "  +0.01%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "58.877.92"
$c.Style = "Normal"

$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("E17").Value = "  -0.74%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.520.13"
$c.Style = "Normal"

$ws.Range("E18").Value = "  +0.98%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "11.09"
$c.Style = "Normal"

$ws.Range("E19").Value = "  +1.23%  "

$ws.Range("E20").Value = "  +0.25%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "323.11"
$c.Style = "Normal"

$ws.Range("E21").Value = "  -0.21%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"

$ws.Range("E22").Value = "  +0.07%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.94"
$c.Style = "Normal"

$ws.Range("E23").Value = "  +2.21%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "65.09"
$c.Style = "Normal"

$ws.Range("E24").Value = "  +3.82%  "

$ws.Range("E25").Value = "  +1.23%  "

$ws.Range("E26").Value = "  -1.45%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"

$ws.Range("E27").Value = "  +0.06%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "7.52"
$c.Style = "Normal"

$ws.Range("E28").Value = "  -1.65%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.0₃0768"
$c.Style = "Normal"

$ws.Range("E29").Value = "  -0.69%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.61"
$c.Style = "Normal"

$ws.Range("E30").Value = "  -0.29%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "171.25"
$c.Style = "Normal"

$ws.Range("E31").Value = "  +3.81%  "

$ws.Range("E32").Value = "  -1.50%  "

$ws.Range("E33").Value = "  +7.99%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"

$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("E35").Value = "  +1.27%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "18.35"
$c.Style = "Normal"

$ws.Range("E36").Value = "  -0.44%  "

$ws.Range("E37").Value = "  -0.78%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.54"
$c.Style = "Normal"

$ws.Range("E38").Value = "  -2.09%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "36.87"
$c.Style = "Normal"

$ws.Range("E39").Value = "  +0.27%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.810"
$c.Style = "Normal"

$ws.Range("E40").Value = "  +1.21%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.59"
$c.Style = "Normal"

$ws.Range("E41").Value = "  -1.08%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "284.84"
$c.Style = "Normal"

$ws.Range("E42").Value = "  +2.59%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.18"
$c.Style = "Normal"

$ws.Range("E43").Value = "  -0.19%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.994"
$c.Style = "Normal"

$ws.Range("E44").Value = "  -0.24%  "

$ws.Range("E45").Value = "  +2.57%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "130.40"
$c.Style = "Normal"

$ws.Range("E46").Value = "  +4.62%  "

$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("E48").Value = "  -1.40%  "

$ws.Range("E49").Value = "  -0.73%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0220"
$c.Style = "Normal"

$ws.Range("E50").Value = "  -0.55%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "17.34"
$c.Style = "Normal"

$ws.Range("E51").Value = "  -1.04%  "
